# Add 'Concentrated quench buffer' to the storage_medium lookup sheet,
# inserted directly after the 'Tris-EDTA' entry, and fix up the
# dependent data validation range on the main sheet.

$wb = $excel.ActiveWorkbook

$storageMedium = $wb.Worksheets.Item("storage_medium")

# Find the row containing "Tris-EDTA" in column A so the new entry is
# inserted right after it (row 11 in the original layout).
$found = $storageMedium.Range("A1:A100").Find("Tris-EDTA")
if ($found -ne $null) {
    $insertRow = $found.Row + 1
} else {
    $insertRow = 12
}

# Insert a new row and shift everything below it down.
$storageMedium.Rows.Item($insertRow).Insert()

$storageMedium.Cells.Item($insertRow, 1).Value = "Concentrated quench buffer"
$storageMedium.Cells.Item($insertRow, 2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000391"

# Update the data validation on the main sheet's storage_medium column (M)
# so the list range covers the new row.
$mainSheet = $wb.Worksheets.Item("Sample Suspension")
$mainSheet.Range("M2:M1001").Validation.Formula1 = "='storage_medium'!`$A`$1:`$A`$23"

# Update the pav:createdOn timestamp on the .metadata sheet.
$metadataSheet = $wb.Worksheets.Item(".metadata")
$metadataSheet.Cells.Item(2, 3).Value = "2024-10-02T11:08:45-07:00"
